$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46074 -> 46075) for every data row (rows 2 through 200).
$ws.Range("C2:C200").Value = 46075
